$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1537058255821364
$ws.Range("D2").Value = 0.05864746209059035
$ws.Range("E2").Value = 0.1087292628355794
$ws.Range("F2").Value = 1.265528993469843
$ws.Range("G2").Value = 1.143207389041763
$ws.Range("H2").Value = 1.120195522848078
$ws.Range("I2").Value = 1.234212713028789
$ws.Range("K2").Value = 1.448946915182546
$ws.Range("L2").Value = 0.1432784980572102
$ws.Range("M2").Value = 0.4299060103059702
$ws.Range("C3").Value = 0.1518217575302572
$ws.Range("D3").Value = 0.05828845029586205
$ws.Range("E3").Value = 0.1087494092808434
$ws.Range("F3").Value = 1.268801857319197
$ws.Range("G3").Value = 1.148937271085273
$ws.Range("H3").Value = 1.129621899936978
$ws.Range("I3").Value = 1.240440210805446
$ws.Range("K3").Value = 1.304713784161095
$ws.Range("L3").Value = 0.1442383484674679
$ws.Range("M3").Value = 0.4015921378362961
$ws.Range("C4").Value = 0.1507170015957655
$ws.Range("D4").Value = 0.05807732462367809
$ws.Range("E4").Value = 0.1088113820550465
$ws.Range("F4").Value = 1.271799893731924
$ws.Range("G4").Value = 1.153443252540413
$ws.Range("H4").Value = 1.136100161980764
$ws.Range("I4").Value = 1.24515460087985
$ws.Range("K4").Value = 1.216251273249554
$ws.Range("L4").Value = 0.1449028153421743
$ws.Range("M4").Value = 0.3843226260682542
$ws.Range("C5").Value = 0.1502799412720037
$ws.Range("D5").Value = 0.05799364817528385
$ws.Range("E5").Value = 0.1088491149013979
$ws.Range("F5").Value = 1.273269550453314
$ws.Range("G5").Value = 1.15552706596408
$ws.Range("H5").Value = 1.13891348436465
$ws.Range("I5").Value = 1.247299249012059
$ws.Range("K5").Value = 1.180227891728435
$ws.Range("L5").Value = 0.145192467651885
$ws.Range("M5").Value = 0.3773144080567548
$ws.Range("C6").Value = 0.150208162453616
$ws.Range("D6").Value = 0.05797989678565685
$ws.Range("E6").Value = 0.1088561342426413
$ws.Range("F6").Value = 1.273528541913052
$ws.Range("G6").Value = 1.155888011854557
$ws.Range("H6").Value = 1.139391099117589
$ws.Range("I6").Value = 1.247668850964779
$ws.Range("K6").Value = 1.174247836415816
$ws.Range("L6").Value = 0.1452417039658123
$ws.Range("M6").Value = 0.3761524728209622
$ws.Range("C7").Value = 0.1507110539996432
$ws.Range("D7").Value = 0.05807618655881441
$ws.Range("E7").Value = 0.1088118404012466
$ws.Range("F7").Value = 1.271818711008663
$ws.Range("G7").Value = 1.153470354168547
$ws.Range("H7").Value = 1.136137401748641
$ws.Range("I7").Value = 1.245182620061463
$ws.Range("K7").Value = 1.215765342990153
$ws.Range("L7").Value = 0.1449066452792245
$ws.Range("M7").Value = 0.3842279920277534
$ws.Range("C8").Value = 0.1530454168845381
$ws.Range("D8").Value = 0.05852175326889508
$ws.Range("E8").Value = 0.1087259172181021
$ws.Range("F8").Value = 1.266451934704186
$ws.Range("G8").Value = 1.144977546651091
$ws.Range("H8").Value = 1.123302296734337
$ws.Range("I8").Value = 1.236174806543978
$ws.Range("K8").Value = 1.399195918769124
$ws.Range("L8").Value = 0.1435938573075006
$ws.Range("M8").Value = 0.4201195911521864
$ws.Range("C9").Value = 0.1580347222860468
$ws.Range("D9").Value = 0.05946866966748843
$ws.Range("E9").Value = 0.1089508875556788
$ws.Range("F9").Value = 1.263800851630123
$ws.Range("G9").Value = 1.136198355918623
$ws.Range("H9").Value = 1.103621963195806
$ws.Range("I9").Value = 1.225600093071016
$ws.Range("K9").Value = 1.759633742078279
$ws.Range("L9").Value = 0.1416160724770137
$ws.Range("M9").Value = 0.4914110271719565
$ws.Range("C10").Value = 0.1619498293260619
$ws.Range("D10").Value = 0.06020816840232612
$ws.Range("E10").Value = 0.1093560819900645
$ws.Range("F10").Value = 1.266696508888771
$ws.Range("G10").Value = 1.134602486366504
$ws.Range("H10").Value = 1.09252506986644
$ws.Range("I10").Value = 1.222185466704666
$ws.Range("K10").Value = 2.024867823284922
$ws.Range("L10").Value = 0.1405275946251727
$ws.Range("M10").Value = 0.544339231139702
$ws.Range("C11").Value = 0.1637848162691284
$ws.Range("D11").Value = 0.0605539454587074
$ws.Range("E11").Value = 0.1095925348392441
$ws.Range("F11").Value = 1.269075088422397
$ws.Range("G11").Value = 1.134942110039702
$ws.Range("H11").Value = 1.088210368319963
$ws.Range("I11").Value = 1.221584775778467
$ws.Range("K11").Value = 2.145617933554888
$ws.Range("L11").Value = 0.1401117927266711
$ws.Range("M11").Value = 0.5685369485545095
$ws.Range("C12").Value = 0.1644874085830423
$ws.Range("D12").Value = 0.06068621673274066
$ws.Range("E12").Value = 0.109689568435428
$ws.Range("F12").Value = 1.270129150674364
$ws.Range("G12").Value = 1.135224862146515
$ws.Range("H12").Value = 1.08668223831414
$ws.Range("I12").Value = 1.221494856183206
$ws.Range("K12").Value = 2.191355518348416
$ws.Range("L12").Value = 0.1399657663520841
$ws.Range("M12").Value = 0.5777171726811332
$ws.Range("C13").Value = 0.1643357499288811
$ws.Range("D13").Value = 0.06065767069532058
$ws.Range("E13").Value = 0.1096683372432068
$ws.Range("F13").Value = 1.269895307810842
$ws.Range("G13").Value = 1.135157096272678
$ws.Range("H13").Value = 1.087006640200599
$ws.Range("I13").Value = 1.221508095504198
$ws.Range("K13").Value = 2.181504592363922
$ws.Range("L13").Value = 0.1399967071574864
$ws.Range("M13").Value = 0.575739291576653
$ws.Range("C14").Value = 0.1638424644192185
$ws.Range("D14").Value = 0.06056480086562743
$ws.Range("E14").Value = 0.1096003676701542
$ws.Range("F14").Value = 1.269158729307762
$ws.Range("G14").Value = 1.134962277716795
$ws.Range("H14").Value = 1.088082526947161
$ws.Range("I14").Value = 1.221574618136295
$ws.Range("K14").Value = 2.149380556513393
$ws.Range("L14").Value = 0.1400995498952256
$ws.Range("M14").Value = 0.5692918702954017
$ws.Range("C15").Value = 0.1635413173207922
$ws.Range("D15").Value = 0.06050808857041545
$ws.Range("E15").Value = 0.1095597102357218
$ws.Range("F15").Value = 1.268727545473595
$ws.Range("G15").Value = 1.134863046940907
$ws.Range("H15").Value = 1.088755319894219
$ws.Range("I15").Value = 1.221633294772609
$ws.Range("K15").Value = 2.12970519834164
$ws.Range("L15").Value = 0.1401640329449094
$ws.Range("M15").Value = 0.565344855042639
$ws.Range("C16").Value = 0.1618309914616134
$ws.Range("D16").Value = 0.06018575840355567
$ws.Range("E16").Value = 0.1093416780202112
$ws.Range("F16").Value = 1.266562481993219
$ws.Range("G16").Value = 1.134601811636202
$ws.Range("H16").Value = 1.092821828012973
$ws.Range("I16").Value = 1.222243940705738
$ws.Range("K16").Value = 2.016978290921031
$ws.Range("L16").Value = 0.1405563665628762
$ws.Range("M16").Value = 0.5427602556110287
$ws.Range("C17").Value = 0.160795559850925
$ws.Range("D17").Value = 0.05999040954750967
$ws.Range("E17").Value = 0.1092212721087513
$ws.Range("F17").Value = 1.265506635391432
$ws.Range("G17").Value = 1.134715125434411
$ws.Range("H17").Value = 1.095504536921183
$ws.Range("I17").Value = 1.222862944554116
$ws.Range("K17").Value = 1.947847037244912
$ws.Range("L17").Value = 0.1408173872948808
$ws.Range("M17").Value = 0.5289360019350084
$ws.Range("C18").Value = 0.1602050906837462
$ws.Range("D18").Value = 0.05987893336882166
$ws.Range("E18").Value = 0.1091569244361246
$ws.Range("F18").Value = 1.264999194217609
$ws.Range("G18").Value = 1.134880544502451
$ws.Range("H18").Value = 1.097116566554433
$ws.Range("I18").Value = 1.223308612543825
$ws.Range("K18").Value = 1.908093475659882
$ws.Range("L18").Value = 0.1409749874247481
$ws.Range("M18").Value = 0.5209960046211251
$ws.Range("C19").Value = 0.1600060423428857
$ws.Range("D19").Value = 0.05984134163395538
$ws.Range("E19").Value = 0.1091359801978022
$ws.Range("F19").Value = 1.26484451120352
$ws.Range("G19").Value = 1.134953741760029
$ws.Range("H19").Value = 1.097674214848709
$ws.Range("I19").Value = 1.22347488563075
$ws.Range("K19").Value = 1.894635174226153
$ws.Range("L19").Value = 0.1410296301947334
$ws.Range("M19").Value = 0.5183096144093611
$ws.Range("C20").Value = 0.1609052574648473
$ws.Range("D20").Value = 0.06001111345840116
$ws.Range("E20").Value = 0.1092335817208578
$ws.Range("F20").Value = 1.265608692289362
$ws.Range("G20").Value = 1.134692681495707
$ws.Range("H20").Value = 1.095211813883353
$ws.Range("I20").Value = 1.222787769759876
$ws.Range("K20").Value = 1.955205265252744
$ws.Range("L20").Value = 0.1407888281617033
$ws.Range("M20").Value = 0.5304064446585102
$ws.Range("C21").Value = 0.1639871450607586
$ws.Range("D21").Value = 0.06059204292549225
$ws.Range("E21").Value = 0.1096201286094889
$ws.Range("F21").Value = 1.269370912702101
$ws.Range("G21").Value = 1.135015309942645
$ws.Range("H21").Value = 1.087763640388758
$ws.Range("I21").Value = 1.221551341310843
$ws.Range("K21").Value = 2.158815845685638
$ws.Range("L21").Value = 0.1400690321718443
$ws.Range("M21").Value = 0.5711851727295425
$ws.Range("C22").Value = 0.1660463316212315
$ws.Range("D22").Value = 0.06097947713367802
$ws.Range("E22").Value = 0.1099164395760042
$ws.Range("F22").Value = 1.272723896503337
$ws.Range("G22").Value = 1.136125060028093
$ws.Range("H22").Value = 1.083512383778782
$ws.Range("I22").Value = 1.221545265998238
$ws.Range("K22").Value = 2.291957281579982
$ws.Range("L22").Value = 0.139665228064672
$ws.Range("M22").Value = 0.5979358400745127
$ws.Range("C23").Value = 0.1649432011377741
$ws.Range("D23").Value = 0.06077199070827533
$ws.Range("E23").Value = 0.1097542966652902
$ws.Range("F23").Value = 1.270852289929024
$ws.Range("G23").Value = 1.135450214701763
$ws.Range("H23").Value = 1.085724844517188
$ws.Range("I23").Value = 1.221474939669648
$ws.Range("K23").Value = 2.220891209666433
$ws.Range("L23").Value = 0.1398746437910461
$ws.Range("M23").Value = 0.5836494873919946
$ws.Range("C24").Value = 0.1608556481763088
$ws.Range("D24").Value = 0.06000175062494151
$ws.Range("E24").Value = 0.1092280013592735
$ws.Range("F24").Value = 1.265562242227972
$ws.Range("G24").Value = 1.134702516114245
$ws.Range("H24").Value = 1.09534393685361
$ws.Range("I24").Value = 1.222821476612744
$ws.Range("K24").Value = 1.951878638572964
$ws.Range("L24").Value = 0.1408017162734225
$ws.Range("M24").Value = 0.5297416333462195
$ws.Range("C25").Value = 0.156641070554798
$ws.Range("D25").Value = 0.05920475668520453
$ws.Range("E25").Value = 0.1088479003459639
$ws.Range("F25").Value = 1.263670648254958
$ws.Range("G25").Value = 1.137724762857928
$ws.Range("H25").Value = 1.108356718106677
$ws.Range("I25").Value = 1.227698477507722
$ws.Range("K25").Value = 1.662050238773872
$ws.Range("L25").Value = 0.1420871539216826
$ws.Range("M25").Value = 0.4720279430649228
